$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.349.04"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'3.811.27"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'700.77"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'171.09"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "'3.810.73"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").Value = "'7.51"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "'0.480"
$ws.Range("E12").Value = "  +4.74%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "'36.00"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'4.453.58"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'3.803.61"
$ws.Range("D17").Value = "'71.355.53"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'7.21"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'17.52"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "'514.95"
$ws.Range("E21").Value = "  +4.10%  "
$ws.Range("D22").Value = "'10.48"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "'84.09"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "'0.0000142"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("D26").Value = "'12.49"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("D27").Value = "'3.958.13"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'10.37"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.35"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'2.23"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").Value = "'29.14"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").Value = "'0.173"
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("D36").Value = "'9.18"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'3.773.75"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.101"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'6.43"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("D41").Value = "'2.42"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "'3.26"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D45").Value = "'171.98"
$ws.Range("E45").Value = "  +5.24%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'0.000308"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").Value = "'49.78"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").Value = "'425.69"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "'8.57"
$ws.Range("E51").Value = "  -0.54%  "
